# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates (row -> new value)
$exhibitUpdates = @{
    "F3"  = 1002
    "F4"  = 782
    "F7"  = 681
    "F8"  = 155
    "F10" = 708
    "F11" = 411
    "F15" = 945
    "F16" = 14
    "F20" = 582
    "F21" = 142
    "F22" = 630
    "F24" = 974
}

foreach ($cell in $exhibitUpdates.Keys) {
    $wsExhibit.Range($cell).Value = $exhibitUpdates[$cell]
}

# 全部类型 sheet updates (row -> new value)
$allUpdates = @{
    "F4"  = 873
    "F5"  = 1002
    "F6"  = 782
    "F9"  = 681
    "F10" = 155
    "F12" = 708
    "F15" = 411
    "F20" = 945
    "F22" = 14
    "F28" = 582
    "F33" = 142
    "F34" = 630
    "F36" = 974
}

foreach ($cell in $allUpdates.Keys) {
    $wsAll.Range($cell).Value = $allUpdates[$cell]
}
